$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test results in the scenario table
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "Failed"
$ws.Range("H2").Value = "30/05/2020"

$ws.Range("B3").Value = "No"
$ws.Range("H3").Value = "26/05/2020"

$ws.Range("B5").Value = "Yes"
$ws.Range("H5").Value = "30/05/2020"

# Update the active cell selection
$ws.Range("B6").Select()
